# "un peu de traduction" - add a new "Vendeurs / Sellers" translation entry
# (MENU.USERS + USERS.TITLE keys) to the translation orderboard sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the two new rows. The order below reproduces the exact order in
# which the new shared strings were first introduced by the author.
$ws.Range("B14").Value = "Vendeurs"
$ws.Range("C14").Value = "Sellers"
$ws.Range("A15").Value = "USERS.TITLE"
$ws.Range("B15").Value = "Liste des vendeurs"
$ws.Range("A14").Value = "MENU.USERS"
$ws.Range("C15").Value = "Sellers"

# Move / leave the active selection where the author left it after editing.
$ws.Range("B18").Select()
